$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values. Values are set with a leading quote-prefix
# marker so Excel treats number-like strings (e.g. "0.998") as literal
# text instead of coercing them to numeric cells, then ClearFormats()
# strips the quote-prefix formatting flag so no stray style index is left
# behind (matches the original un-styled inline-string cells).

$ws.Range("D2").Value = '''63.665.78'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  -0.50%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''2.723.39'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  -1.22%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''561.25'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -2.35%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''158.04'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  -0.58%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '''  -0.02%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '''  -1.49%  '
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = '''  -2.22%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = '''0.166'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '''  +0.34%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''5.59'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  -3.12%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = '''  -2.90%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''3.204.52'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  -1.28%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = '''26.66'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  -0.94%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = '''63.533.01'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '''  -0.11%  '
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = '''  -2.55%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''2.724.90'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  -1.39%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = '''12.19'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  +0.42%  '
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = '''  -3.81%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = '''349.35'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  -1.90%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = '''  -3.74%  '
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = '''  +0.19%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  -2.19%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = '''64.15'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '''  -1.15%  '
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = '''  +0.17%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''0.998'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  -0.10%  '
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = '''  -3.96%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = '''0.0₃0886'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '''  -1.59%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D30").Value = '''1.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  +0.36%  '
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = '''  -1.73%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = '''164.87'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '''  -2.69%  '
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = '''EthereumClassic'
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = '''19.88'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '''  -1.12%  '
$ws.Range("E33").ClearFormats()
$ws.Range("B34").Value = '''NEARProtocol'
$ws.Range("B34").ClearFormats()
$ws.Range("C34").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C34").ClearFormats()
$ws.Range("D34").Value = '''4.85'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '''  -1.01%  '
$ws.Range("E34").ClearFormats()
$ws.Range("B35").Value = '''USDe'
$ws.Range("B35").ClearFormats()
$ws.Range("C35").Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C35").ClearFormats()
$ws.Range("D35").Value = '''0.998'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  -0.02%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = '''  -0.30%  '
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = '''  +0.17%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''348.83'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  -0.41%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = '''0.961'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  -4.48%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = '''6.16'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '''  -1.35%  '
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = '''  -3.95%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = '''38.32'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '''  -1.96%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''21.43'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  -1.40%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''20.72'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  -3.33%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = '''0.0575'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  -2.41%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = '''0.628'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  -0.81%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = '''132.54'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '''  -3.24%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = '''0.998'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  -0.06%  '
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = '''  -2.80%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''11.08'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  +0.26%  '
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = '''  -2.95%  '
$ws.Range("E51").ClearFormats()
